# Apply the "Fixing network data cleaning scripts" edit to NEW_HAMPSHIRE_2019
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns to snake_case machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Normalize capitalization of connector words ("de"/"el" -> "De"/"El")
#    in place names throughout the data rows
$ws.Range("B11").Value = "San Cristóbal De Las Casas"
$ws.Range("A19").Value = "Ciudad De México"
$ws.Range("A26").Value = "Estado De México"
$ws.Range("B32").Value = "Apaseo El Alto"
$ws.Range("B38").Value = "Mártir De Cuilapan"
$ws.Range("B39").Value = "Taxco De Alarcón"
$ws.Range("B51").Value = "Ocotlán De Morelos"
$ws.Range("B54").Value = "Huehuetlán El Chico"
$ws.Range("B59").Value = "Jalpan De Serra"
$ws.Range("B60").Value = "Landa De Matamoros"
$ws.Range("B69").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B72").Value = "Soledad De Doblado"

# 3) Remove the trailing metadata/footer rows (78-82): sample size, source,
#    author and date notes that don't belong in the cleaned dataset
$ws.Rows("78:82").Delete()
